$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.458.63"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "'1.624.68"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'211.78"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'0.0860"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "'1.857.92"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'1.633.95"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'0.557"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "'27.421.54"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "'228.52"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'0.0₃0719"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'7.51"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'10.66"
$ws.Range("E22").Value = "  +6.48%  "
$ws.Range("D23").Value = "'4.35"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").Value = "'149.47"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'6.87"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.111"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "'15.55"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'0.0480"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "'1.467.02"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").Value = "'2.33"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.558"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.921"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0167"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.873"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").Value = "'1.01"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "'1.02"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'67.56"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'2.27"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").Value = "'5.36"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("D47").Value = "'1.764.29"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").Value = "'87.23"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0987"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.0₇0992"
$ws.Range("E51").Value = "  -5.96%  "
